# Reward_Table.xlsx edit: add a new "melon" (🍈) slot symbol, re-balance the
# probability column, replace the Lemon-combo reward row with a Melon-combo
# reward row, and halve several of the lower-tier reward payouts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- A/B columns: symbol + probability table (A1:B7) -----------------------
# Row 7 is brand new (didn't exist before); B7 used to hold the `=SUM(B1:B6)`
# formula, which now lives in B8 instead.
$ws.Range("A1").Value = "🍋"
$ws.Range("B1").Value = 0.25

$ws.Range("A2").Value = "🍒"
$ws.Range("B2").Value = 0.225

$ws.Range("A3").Value = "🍊"
$ws.Range("B3").Value = 0.2

$ws.Range("A4").Value = "🍈"
$ws.Range("B4").Value = 0.15

$ws.Range("A5").Value = "🍇"
$ws.Range("B5").Value = 0.1

$ws.Range("A6").Value = "🍉"
$ws.Range("B6").Value = 0.05

$ws.Range("A7").Value = "💎"
$ws.Range("B7").Value = 0.025

# --- Reward table rows 2-6: labels/amounts unchanged, nothing to do here ---

# --- Row 7 becomes the new Melon (🍈🍈🍈) reward row -------------------------
$ws.Range("D7").Value = "🍈🍈🍈"
$ws.Range("E7").Value = 50
$ws.Range("F7").Value = "🍈🍈🍈"
$ws.Range("G7").Value = 100
$ws.Range("H7").Value = "🍈🍈🍈"
$ws.Range("I7").Value = 150
$ws.Range("J7").Value = "🍈🍈🍈"
# The old J7 carried the bold/black-font style (s="2"); the fresh melon row
# doesn't, so strip any inherited formatting back to the default.
$ws.Range("J7").ClearFormats()

# --- Row 8: SUM formula (shifted down from B7) + Melon (🍈🍈💎) reward row --
$ws.Range("B8").Formula = "=SUM(B1:B7)"

$ws.Range("D8").Value = "🍈🍈💎"
$ws.Range("E8").Value = 50
$ws.Range("F8").Value = "🍈🍈💎"
$ws.Range("G8").Value = 100
$ws.Range("H8").Value = "🍈🍈💎"
$ws.Range("I8").Value = 150
$ws.Range("J8").Value = "🍈🍈💎"
$ws.Range("J8").ClearFormats()

# --- Rows 9-13: labels unchanged, payouts halved ----------------------------
$ws.Range("E9").Value = 30
$ws.Range("G9").Value = 60
$ws.Range("I9").Value = 90

$ws.Range("E10").Value = 30
$ws.Range("G10").Value = 60
$ws.Range("I10").Value = 90

$ws.Range("E11").Value = 10
$ws.Range("G11").Value = 20
$ws.Range("I11").Value = 30

$ws.Range("E12").Value = 5
$ws.Range("G12").Value = 10
$ws.Range("I12").Value = 15

$ws.Range("E13").Value = 1
$ws.Range("G13").Value = 2
$ws.Range("I13").Value = 3

# --- Selection moves to B12 --------------------------------------------------
$ws.Range("B12").Select()
